# Updates crypto price/volume figures (and restores the ShibaInu /
# WrappedliquidstakedEther2.0 row ordering) per the latest GitHub Actions
# refresh of the cryptos list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '73.871.21'
$ws.Range('E2').Value = '  +7.45%  '
# Row 3
$ws.Range('D3').Value = '2.627.08'
$ws.Range('E3').Value = '  +7.80%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '184.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +14.54%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '582.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.32%  '
# Row 7
$ws.Range('E7').Value = '  -0.16%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.534'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.26%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.201'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +19.62%  '
# Row 10
$ws.Range('D10').Value = '2.625.76'
$ws.Range('E10').Value = '  +7.78%  '
# Row 11
$ws.Range('E11').Value = '  +0.28%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.358'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.32%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.77'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.11%  '
# Row 14
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000189'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.91%  '
# Row 15
$ws.Range('D15').Value = '73.701.08'
$ws.Range('E15').Value = '  +7.36%  '
# Row 16
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '3.102.66'
$ws.Range('E16').Value = '  +7.54%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.18'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +12.88%  '
# Row 18
$ws.Range('D18').Value = '2.621.93'
$ws.Range('E18').Value = '  +7.58%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.10'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +31.30%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.88'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +12.52%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '371.62'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +9.56%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +19.17%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.98%  '
# Row 24
$ws.Range('E24').Value = '  +0.05%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '69.81'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.30%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.15'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.86%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.42'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +14.71%  '
# Row 28
$ws.Range('D28').Value = '2.748.13'
$ws.Range('E28').Value = '  +7.20%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.75%  '
# Row 30
$ws.Range('D30').Value = '0.0₃0942'
$ws.Range('E30').Value = '  +15.11%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '523.42'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +22.35%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.40'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +21.03%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.66'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.50%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.74'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.14%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.07%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.120'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +13.16%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '160.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.87%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.16'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.58%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.26'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.43%  '
# Row 40
$ws.Range('E40').Value = '  -0.06%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.91'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +13.05%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.328'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.87%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.67'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +11.01%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '161.49'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +23.30%  '
# Row 45
$ws.Range('E45').Value = '  +10.42%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.37'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +16.01%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '38.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.16%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0851'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +18.61%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.63'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +8.97%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.528'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +10.01%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.72'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +22.79%  '
